$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.279.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.507"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.05"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.874.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.622.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.55"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.241.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "196.06"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.53"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.98"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.64"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0509"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.139.49"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.556"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.50"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.30"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.782.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.35"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.419"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.71"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.14%  "
